$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 18 de Julio de 2020 a las 06:19"

# Row 15 - Pakistan
$ws.Range("B15").Value = 261916
$ws.Range("C15").Value = 1917
$ws.Range("D15").Value = 198509
$ws.Range("E15").Value = 57885
$ws.Range("G15").Value = 47
$ws.Range("H15").Value = 5522

# Row 32 - Kazajistan
$ws.Range("B32").Value = 68703
$ws.Range("C32").Value = 1808
$ws.Range("E32").Value = 28072

# Row 35 - Belgica
$ws.Range("B35").Value = 63499
$ws.Range("C35").Value = 261
$ws.Range("D35").Value = 17289
$ws.Range("E35").Value = 36410
$ws.Range("G35").Value = 5
$ws.Range("H35").Value = 9800

# Row 55 - Honduras
$ws.Range("B55").Value = 31745
$ws.Range("C55").Value = 878
$ws.Range("D55").Value = 3565
$ws.Range("E55").Value = 27323
$ws.Range("G55").Value = 22
$ws.Range("H55").Value = 857

# Row 89 - Haiti
$ws.Range("B89").Value = 6975
$ws.Range("C89").Value = 27
$ws.Range("D89").Value = 3738
$ws.Range("E89").Value = 3091
$ws.Range("G89").Value = 1
$ws.Range("H89").Value = 146

# Row 170 - Mongolia
$ws.Range("B170").Value = 287
$ws.Range("C170").Value = 25
$ws.Range("E170").Value = 76
